$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ependorf 200, 50, 20 uL")

$ws.Range("B2").Value = 998.4
$ws.Range("E2").Value = 998.2
$ws.Range("H2").Value = 997.7
$ws.Range("B3").Value = 1000
$ws.Range("E3").Value = 997.2
$ws.Range("H3").Value = 998.2
$ws.Range("B4").Value = 998.4
$ws.Range("E4").Value = 999.6
$ws.Range("H4").Value = 997.2
$ws.Range("B5").Value = 1004
$ws.Range("E5").Value = 1000
$ws.Range("H5").Value = 996.8
$ws.Range("B6").Value = 1004
$ws.Range("E6").Value = 1000.3
$ws.Range("H6").Value = 997.7
$ws.Range("B10").Value = 502.7
$ws.Range("E10").Value = 500
$ws.Range("H10").Value = 501
$ws.Range("B11").Value = 502.3
$ws.Range("E11").Value = 501.8
$ws.Range("H11").Value = 501.6
$ws.Range("B12").Value = 502.1
$ws.Range("E12").Value = 501.4
$ws.Range("H12").Value = 500.5
$ws.Range("B13").Value = 502.4
$ws.Range("E13").Value = 500.5
$ws.Range("H13").Value = 501
$ws.Range("B14").Value = 502
$ws.Range("E14").Value = 502.1
$ws.Range("H14").Value = 501.2
$ws.Range("B18").Value = 203.4
$ws.Range("E18").Value = 203.7
$ws.Range("H18").Value = 203
$ws.Range("B19").Value = 203.5
$ws.Range("E19").Value = 203.7
$ws.Range("H19").Value = 203.8
$ws.Range("B20").Value = 203
$ws.Range("E20").Value = 203.8
$ws.Range("H20").Value = 203.6
$ws.Range("E21").Value = 203.5
$ws.Range("H21").Value = 203
$ws.Range("B22").Value = 203.7
$ws.Range("E22").Value = 203.8
$ws.Range("H22").Value = 204
$ws.Range("B26").Value = 104.2
$ws.Range("E26").Value = 105
$ws.Range("H26").Value = 104.1
$ws.Range("B27").Value = 104.5
$ws.Range("E27").Value = 105.6
$ws.Range("H27").Value = 104.1
$ws.Range("B28").Value = 104.6
$ws.Range("E28").Value = 105.4
$ws.Range("H28").Value = 104.1
$ws.Range("B29").Value = 105
$ws.Range("E29").Value = 105.1
$ws.Range("H29").Value = 103.7
$ws.Range("B30").Value = 104.9
$ws.Range("E30").Value = 104.8
$ws.Range("H30").Value = 102
